# Add the August 2025 actuals column (K) to the historical data table.
# Previously the system only carried a hardcoded 7-month window (through
# column J, which just duplicated July's column I). Now that the bug is
# fixed, the next available month of real data (August) is appended as
# its own column so the September forecast can use it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 9234.4775200000004
$ws.Range("K3").Value = 32092.500000000004
$ws.Range("K4").Value = 14196.32
$ws.Range("K5").Value = 25595.16
$ws.Range("K6").Value = 5004.93
